$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the trailing blank spacer row (old row 12); old row 13 becomes new row 12.
$ws.Rows("12").Delete()

# Insert a new column before the old "COMMUNITY" column (old E, becomes F) to hold
# the new "RELEASE DATE" column.
$ws.Columns("E").Insert()
$ws.Columns("E").ColumnWidth = 19.16666666666667

# Header
$ws.Range("E3").Value = "RELEASE DATE"

# Release date values (Excel serial dates)
$ws.Range("E4").Value = 40316
$ws.Range("E5").Value = 38762
$ws.Range("E6").Value = 37652
$ws.Range("E7").Value = 40803
$ws.Range("E8").Value = 39769
$ws.Range("E9").Value = 43304
$ws.Range("E10").Value = 40469
$ws.Range("E11").Value = 43678

# Number formats per row
$ws.Range("E4:E6").NumberFormat = "mmm-yy"
$ws.Range("E7:E8").NumberFormat = "d-mmm-yy"
$ws.Range("E9").NumberFormat = "mm-dd-yy"
$ws.Range("E10:E11").NumberFormat = "d-mmm-yy"

# Match the author's final selection
$ws.Range("A12").Select()
